$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1499.8
$ws.Range("I40").Value = 1499.8334
$ws.Range("J40").Value = 1499.75
$ws.Range("K40").Value = 1499.8334
$ws.Range("L40").Value = 1499.75
$ws.Range("M40").Value = -1324.8334
$ws.Range("N40").Value = -1849.75

$ws.Range("H64").Value = 2917.5881
$ws.Range("J64").Value = 2999.9
$ws.Range("L64").Value = 2999.9
$ws.Range("N64").Value = -3495.9

$ws.Range("H67").Value = 2917.5881
$ws.Range("J67").Value = 2999.9
$ws.Range("L67").Value = 2999.9
$ws.Range("N67").Value = -4715.9

$ws.Range("H76").Value = 4904920
$ws.Range("I76").Value = 5750078.5
$ws.Range("K76").Value = 5750078.5
$ws.Range("M76").Value = -5749763.5

$ws.Range("H79").Value = 4904920
$ws.Range("I79").Value = 5750078.5
$ws.Range("K79").Value = 5750078.5
$ws.Range("M79").Value = -5748986.5

$ws.Range("H112").Value = 11848.125
$ws.Range("J112").Value = 11848.125
$ws.Range("L112").Value = 35544.375
$ws.Range("N112").Value = -37760.375

$ws.Range("H129").Value = 1304.0513
$ws.Range("J129").Value = 1794.12
$ws.Range("L129").Value = 5382.36
$ws.Range("N129").Value = -15382.36

$ws.Range("H132").Value = 4555.6
$ws.Range("I132").Value = 4434.0557
$ws.Range("J132").Value = 5649.5
$ws.Range("K132").Value = 13302.1671
$ws.Range("L132").Value = 16948.5
$ws.Range("M132").Value = -10772.1671
$ws.Range("N132").Value = -22008.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 36637.5
$ws.Range("I63").Value = 104125
$ws.Range("J63").Value = 4669.737
$ws.Range("K63").Value = 104125
$ws.Range("L63").Value = 4669.737
$ws.Range("M63").Value = -103439
$ws.Range("N63").Value = -6041.737

$ws.Range("H66").Value = 36637.5
$ws.Range("I66").Value = 104125
$ws.Range("J66").Value = 4669.737
$ws.Range("K66").Value = 520625
$ws.Range("L66").Value = 23348.685
$ws.Range("M66").Value = -517193
$ws.Range("N66").Value = -30212.685

$ws.Range("H88").Value = 2486.5715
$ws.Range("I88").Value = 2035.3334
$ws.Range("J88").Value = 2825
$ws.Range("K88").Value = 2035.3334
$ws.Range("L88").Value = 2825
$ws.Range("M88").Value = -1629.3334
$ws.Range("N88").Value = -3637

$ws.Range("H91").Value = 2486.5715
$ws.Range("I91").Value = 2035.3334
$ws.Range("J91").Value = 2825
$ws.Range("K91").Value = 2035.3334
$ws.Range("L91").Value = 2825
$ws.Range("M91").Value = -631.3334
$ws.Range("N91").Value = -5633

$ws.Range("H129").Value = 49996
$ws.Range("J129").Value = 49996
$ws.Range("L129").Value = 49996
$ws.Range("N129").Value = -59996

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 4000
$ws.Range("I105").Value = 4000
$ws.Range("K105").Value = 4000
$ws.Range("M105").Value = -2253

$ws.Range("H141").Value = 20000
$ws.Range("J141").Value = 20000
$ws.Range("L141").Value = 20000
$ws.Range("N141").Value = -30360

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2003.6154
$ws.Range("I16").Value = 2079.125
$ws.Range("J16").Value = 1882.8
$ws.Range("K16").Value = 2079.125
$ws.Range("L16").Value = 1882.8
$ws.Range("M16").Value = -1792.125
$ws.Range("N16").Value = -2456.8

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("N41").ClearContents()

$ws.Range("H62").Value = 4500
$ws.Range("I62").Value = 4500
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 4500
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -3876
$ws.Range("N62").ClearContents()

$ws.Range("H65").Value = 4500
$ws.Range("I65").Value = 4500
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 22500
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -19380
$ws.Range("N65").ClearContents()

$ws.Range("H113").Value = 2003.6154
$ws.Range("I113").Value = 2079.125
$ws.Range("J113").Value = 1882.8
$ws.Range("K113").Value = 2079.125
$ws.Range("L113").Value = 1882.8
$ws.Range("M113").Value = 90.875
$ws.Range("N113").Value = -6222.8

$ws.Range("H122").Value = 1922.8695
$ws.Range("I122").Value = 1804
$ws.Range("J122").Value = 1999.2858
$ws.Range("K122").Value = 5412
$ws.Range("L122").Value = 5997.857400000001
$ws.Range("M122").Value = -2962
$ws.Range("N122").Value = -10897.8574

$ws.Range("H133").Value = 48850
$ws.Range("J133").Value = 48850
$ws.Range("L133").Value = 48850
$ws.Range("N133").Value = -53910

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4119665.5
$ws.Range("J4").Value = 2536.04
$ws.Range("L4").Value = 7608.12
$ws.Range("N4").Value = -7832.12

$ws.Range("H5").Value = 710.36365
$ws.Range("I5").Value = 487.63635
$ws.Range("J5").Value = 1044.4546
$ws.Range("K5").Value = 1462.90905
$ws.Range("L5").Value = 3133.3638
$ws.Range("M5").Value = -1350.90905
$ws.Range("N5").Value = -3357.3638

$ws.Range("H6").Value = 321.54544
$ws.Range("I6").Value = 31.416666
$ws.Range("J6").Value = 669.7
$ws.Range("K6").Value = 94.24999800000001
$ws.Range("L6").Value = 2009.1
$ws.Range("M6").Value = 18.75000199999999
$ws.Range("N6").Value = -2235.1

$ws.Range("H25").Value = 485
$ws.Range("J25").Value = 500
$ws.Range("L25").Value = 1500
$ws.Range("N25").Value = -1838

$ws.Range("H30").Value = 485
$ws.Range("J30").Value = 500
$ws.Range("L30").Value = 1500
$ws.Range("N30").Value = -1704

$ws.Range("H135").Value = 710.36365
$ws.Range("I135").Value = 487.63635
$ws.Range("J135").Value = 1044.4546
$ws.Range("K135").Value = 4388.72715
$ws.Range("L135").Value = 9400.091400000001
$ws.Range("M135").Value = -1853.72715
$ws.Range("N135").Value = -14470.0914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10662.117
$ws.Range("I70").Value = 12688
$ws.Range("K70").Value = 12688
$ws.Range("M70").Value = -12418

$ws.Range("H73").Value = 10662.117
$ws.Range("I73").Value = 12688
$ws.Range("K73").Value = 12688
$ws.Range("M73").Value = -11752

$ws.Range("H80").Value = 1568519.5
$ws.Range("I80").Value = 3001701.8
$ws.Range("J80").Value = 135337.33
$ws.Range("K80").Value = 3001701.8
$ws.Range("L80").Value = 135337.33
$ws.Range("M80").Value = -3000703.8
$ws.Range("N80").Value = -137333.33

$ws.Range("H83").Value = 1568519.5
$ws.Range("I83").Value = 3001701.8
$ws.Range("J83").Value = 135337.33
$ws.Range("K83").Value = 15008509
$ws.Range("L83").Value = 676686.6499999999
$ws.Range("M83").Value = -15003517
$ws.Range("N83").Value = -686670.6499999999

$ws.Range("H102").Value = 1498.125
$ws.Range("I102").Value = 1495.8572
$ws.Range("K102").Value = 1495.8572
$ws.Range("M102").Value = 126.1428000000001

$ws.Range("H140").Value = 70779.75
$ws.Range("J140").Value = 70779.75
$ws.Range("L140").Value = 70779.75
$ws.Range("N140").Value = -81139.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2123.5454
$ws.Range("I40").Value = 1928.7778
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 1928.7778
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1792.7778
$ws.Range("N40").Value = -3272

$ws.Range("H61").Value = 2730.4285
$ws.Range("I61").Value = 2193.5652
$ws.Range("J61").Value = 5200
$ws.Range("K61").Value = 2193.5652
$ws.Range("L61").Value = 5200
$ws.Range("M61").Value = -1991.5652
$ws.Range("N61").Value = -5604

$ws.Range("H113").Value = 2730.4285
$ws.Range("I113").Value = 2193.5652
$ws.Range("J113").Value = 5200
$ws.Range("K113").Value = 2193.5652
$ws.Range("L113").Value = 5200
$ws.Range("M113").Value = -23.5652
$ws.Range("N113").Value = -9540

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 1116.5
$ws.Range("I100").Value = 1980
$ws.Range("J100").Value = 828.6667
$ws.Range("K100").Value = 3960
$ws.Range("L100").Value = 1657.3334
$ws.Range("M100").Value = -3419
$ws.Range("N100").Value = -2739.3334
